$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = 10.0
$ws.Range("F6").Value = "Sprengung"
$ws.Range("G6").Value = 750.0
